$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Add the new "Sheet1" worksheet (tab order: after tt_simulation_1)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet1"

# Fill in the win/loss analysis table on the new sheet.
$ws2.Range("A2").Value = 11
$ws2.Range("B2").Value = 5
$ws2.Range("C2").Value = "w"
$ws2.Range("D2").Formula = '=ABS(B2-A2)/11'

$ws2.Range("A3").Value = 11
$ws2.Range("B3").Value = 9
$ws2.Range("C3").Value = "w"
$ws2.Range("D3:D6").Formula = '=ABS(B3-A3)/11'

$ws2.Range("A4").Value = 9
$ws2.Range("B4").Value = 11
$ws2.Range("C4").Value = "l"

$ws2.Range("A5").Value = 8
$ws2.Range("B5").Value = 11
$ws2.Range("C5").Value = "l"

$ws2.Range("A6").Value = 11
$ws2.Range("B6").Value = 9
$ws2.Range("C6").Value = "w"

$ws2.Range("C8").Value = "winner_average"
$ws2.Range("D8").Formula = '=AVERAGEIF(C2:C6,"w",D2:D6)'
$ws2.Range("E8").Value = 3
$ws2.Range("F8").Formula = '=D8*3+E8'

$ws2.Range("C9").Value = "loser_average"
$ws2.Range("D9").Formula = '=AVERAGEIF(C3:C7,"l",D3:D7)'
$ws2.Range("E9").Value = 2
$ws2.Range("F9").Formula = '=D9*2+E9'

$ws2.Range("F13").Select()

# ---------------------------------------------------------------------------
# 2. Rework the scoring formulas on tt_simulation_1 (sheet 1)
# ---------------------------------------------------------------------------

# Column D: margin = 5 - wins-against (C)
$ws1.Range("D2").Formula = '=5-C2'
$ws1.Range("D3:D13").Formula = '=5-C3'

# Column G: total points now also credits the bonus points earned from wins
$ws1.Range("G2").Formula = '=E2+C2*F2'
$ws1.Range("G3:G13").Formula = '=E3+C3*F3'

# Column I: new ranking "delta" split out from the previous H formula.
# Clone the formatting from column H first so I picks up the same cell
# style (border + centered alignment) as the rest of the table.
$ws1.Range("H2").Copy()
$ws1.Range("I2:I13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws1.Range("I2").Formula = '=12-ROUND(12*(G2/$H$16),1)'
$ws1.Range("I3:I13").Formula = '=12-ROUND(12*(G3/$H$16),1)'

# H16 used to be MAX(G2:G13); it becomes a fixed reference value instead.
$ws1.Range("G16").ClearContents()
$ws1.Range("F16").Value = "max total points"
$ws1.Range("F16:G16").HorizontalAlignment = -4108
$ws1.Range("F16:G16").Merge()
$ws1.Range("H16").Formula = '10'

$ws1.Range("J8").Select()
